# Add a "Notes" sheet (changelog) as the last tab, after "Dates".
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Add()
$ws.Name = "Notes"

# Fill in the changelog content first (while this is still the active/fresh
# sheet reference) - moving the sheet before writing can leave the writes
# targeting the wrong tab.
$ws.Range("B2").Value = "Notes"

$ws.Range("B3").Value = 1
$ws.Range("C3").Value = "Added data for 2016"

$ws.Range("B4").Value = 2
$ws.Range("C4").Value = "In 2016 tab, added Delivered By column"

$ws.Range("B5").Value = 3
$ws.Range("C5").Value = "Summarized 2016 data in Dates tab"

$ws.Range("B6").Value = 4
$ws.Range("C6").Value = "Added data for 2017"

$ws.Range("B7").Value = "…"

# Move the new sheet to the end of the tab strip, right after "Dates".
$datesSheet = $wb.Worksheets.Item("Dates")
$ws.Move($null, $datesSheet)
